# LMSTestData.xlsx batch edit / integration changes
# Target sheet: "BatchPage" (2nd worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Update text of the two surviving rows that changed value (do this BEFORE
#    deleting rows, and in this order, so the shared-string table ends up
#    with the same ordering as the target workbook).
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Team7-TechTesters@#-SDET-SDET01-01123"
$ws.Range("B2").Value = "Team7-TechTesters-SDET-SDET01-09"

# ---------------------------------------------------------------------------
# 2. Remove the rows that were dropped from the test-data grid. Deleting from
#    the bottom up keeps the remaining row numbers stable while we work.
#    Old row 13 : duplicate "Batchstatus is null" scenario
#    Old row 7  : "BatchStatus is Inactive" scenario
#    Old row 6  : "Empty BatchNoOfclasses" scenario
#    Old row 3  : duplicate "Successfully Created" scenario (lower-case name)
# ---------------------------------------------------------------------------
$ws.Range("A13").EntireRow.Delete()
$ws.Range("A7").EntireRow.Delete()
$ws.Range("A6").EntireRow.Delete()
$ws.Range("A3").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3. Re-point the hyperlink that used to live on B5 (old layout) so it now
#    points at B4 (new layout) with the updated mailto target. Adding a
#    hyperlink auto-applies the "Hyperlink" cell style, which the source
#    cell never had, so explicitly restore the plain/Normal style and drop
#    the now-unused "Hyperlink" named cell style again.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Team7-TechTesters@#-SDET-SDET01-01123")
$ws.Range("B4").Style = "Normal"
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 4. Update the active selection to match the new layout.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B2").Select()
